$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts old D:K data to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number formatting/style from column F (the old column D, now shifted) into the new D:E columns
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)

# Populate the two new quarter columns (D = 2018-12-31, E = 2018-09-30) with their reported values
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 112400
$ws.Cells.Item(8, 5).Value = 101400
$ws.Cells.Item(9, 4).Value = 6400
$ws.Cells.Item(9, 5).Value = 5900
$ws.Cells.Item(10, 4).Value = 106000
$ws.Cells.Item(10, 5).Value = 95500
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 5800
$ws.Cells.Item(15, 5).Value = 6200
$ws.Cells.Item(17, 4).Value = 58500
$ws.Cells.Item(17, 5).Value = 55000
$ws.Cells.Item(18, 4).Value = 53900
$ws.Cells.Item(18, 5).Value = 46400
$ws.Cells.Item(20, 4).Value = 2100
$ws.Cells.Item(20, 5).Value = 1400
$ws.Cells.Item(21, 4).Value = 61900
$ws.Cells.Item(21, 5).Value = 54000
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 56100
$ws.Cells.Item(23, 5).Value = 47800
$ws.Cells.Item(24, 4).Value = 10600
$ws.Cells.Item(24, 5).Value = 9200
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 45400
$ws.Cells.Item(26, 5).Value = 38600
$ws.Cells.Item(27, 4).Value = 45400
$ws.Cells.Item(27, 5).Value = 38600
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 400
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -2100
$ws.Cells.Item(32, 5).Value = -1400
$ws.Cells.Item(33, 4).Value = 45800
$ws.Cells.Item(33, 5).Value = 38600
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 45800
$ws.Cells.Item(35, 5).Value = 38600
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 246300
$ws.Cells.Item(41, 5).Value = 210100
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 57500
$ws.Cells.Item(43, 5).Value = 62000
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(45, 4).Value = 22500
$ws.Cells.Item(45, 5).Value = 18600
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(47, 4).Value = 240100
$ws.Cells.Item(47, 5).Value = 236100
$ws.Cells.Item(48, 4).Value = 63000
$ws.Cells.Item(48, 5).Value = 54300
$ws.Cells.Item(49, 4).Value = 62700
$ws.Cells.Item(49, 5).Value = 62800
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 3400
$ws.Cells.Item(52, 5).Value = 2300
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 695500
$ws.Cells.Item(54, 5).Value = 646100
$ws.Cells.Item(57, 4).Value = 29400
$ws.Cells.Item(57, 5).Value = 19300
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(59, 4).Value = 55500
$ws.Cells.Item(59, 5).Value = 44400
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 87700
$ws.Cells.Item(66, 5).Value = 67500
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 463300
$ws.Cells.Item(72, 5).Value = 433200
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 607900
$ws.Cells.Item(76, 5).Value = 578600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 45800
$ws.Cells.Item(81, 5).Value = 38600
$ws.Cells.Item(83, 4).Value = 5800
$ws.Cells.Item(83, 5).Value = 6200
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 70100
$ws.Cells.Item(89, 5).Value = 68000
$ws.Cells.Item(91, 4).Value = -12100
$ws.Cells.Item(91, 5).Value = -7800
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -13300
$ws.Cells.Item(94, 5).Value = -7600
$ws.Cells.Item(96, 4).Value = -15500
$ws.Cells.Item(96, 5).Value = -15500
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -20200
$ws.Cells.Item(100, 5).Value = -20600
$ws.Cells.Item(101, 4).Value = -300
$ws.Cells.Item(101, 5).Value = -500
$ws.Cells.Item(102, 4).Value = 36200
$ws.Cells.Item(102, 5).Value = 39400
